$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Remove the "Required offshore wind" worksheet entirely.
# ------------------------------------------------------------------
$reqWind = $wb.Worksheets.Item("Required offshore wind")
[void]$reqWind.Delete()

# ------------------------------------------------------------------
# 2. Rebuild the "About" sheet: clear everything (contents + formatting)
#    and re-enter the simplified content.
# ------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
[void]$about.Cells.Clear()

$about.Range("A1").Value = "PMCCS Policy Mandated Capacity Construction Schedule"
$about.Range("A1").Font.Bold = $true

$about.Range("A3").Value = "Source:"
$about.Range("A3").Font.Bold = $true
$about.Range("B3").Value = "None (this variable is intended to be user-specified)"

$about.Range("A5").Value = "Note:"
$about.Range("A5").Font.Bold = $true
$about.Range("A6").Value = "You may use this variable to specify the electricity generating capacity"
$about.Range("A7").Value = "in MW that will be built each year in the policy case when the"
$about.Range("A8").Value = """Boolean Use Non BAU Mandated Capacity Construction Schedule"""
$about.Range("A9").Value = "policy is enabled."

# ------------------------------------------------------------------
# 3. PMCCS sheet: drop the array formula driving the "offshore wind"
#    row (row 14) now that its source sheet is gone, resetting the
#    construction values for 2024-2035 (cols I:T) back to 0.
# ------------------------------------------------------------------
$pmccs = $wb.Worksheets.Item("PMCCS")
$pmccs.Range("B14:T14").Value = 0

# ------------------------------------------------------------------
# 4. Fix up sheet selections / active tab to match the saved state:
#    "About" becomes the active tab, PMCCS keeps a plain selection.
# ------------------------------------------------------------------
$pmccs.Range("A18").Select()
$about.Select()
$about.Range("A1").Select()
